$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "shopping with Botticelli"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1983-05-02"

$ws.Range("E2").Value = "Intelligent Wool Chair"
$ws.Range("F2").Value = "Iceberg lettuce"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2.0"
